$wb = $excel.ActiveWorkbook

# Sheet "OFF" - update row 3 (R row) values
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B3").Value = 500
$wsOff.Range("C3").Value = 346
$wsOff.Range("D3").Value = 104
$wsOff.Range("E3").Value = 55
$wsOff.Range("F3").Value = 12

# Sheet "DEF" - update row 3 (R row) values
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B3").Value = 405
$wsDef.Range("C3").Value = 281
$wsDef.Range("D3").Value = 114
$wsDef.Range("E3").Value = 66
